$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The quarterly report gained a new "12 ماهه منتهی به 1401/12" column.
# Every existing reporting column shifts one slot to the left (the
# oldest period, "6 ماهه منتهی به 1399/06", drops off) and the brand
# new period is appended as the new column M.
# ------------------------------------------------------------------

# Step 1: drop the oldest reporting column and shift everything left
$ws.Range("D:D").Delete()

# Step 1b: carry the blank-row formatting into the newly exposed column M
# (these rows have no data, just the row's background/border style)
$ws.Range("L1:L7").Copy()
$ws.Range("M1:M7").PasteSpecial(-4122)
$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial(-4122)
$ws.Range("L28").Copy()
$ws.Range("M28").PasteSpecial(-4122)

# Step 1c: column M now holds an annual ("12 ماهه") period like E and I,
# so it gets the same wider column formatting used for the other annual columns
$ws.Range("M1").ColumnWidth = 28.15

# Step 2: label the new column (period header + publish-date header)
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-23 (2)"

# Step 3: the 9-month-1400/09 report got a new revision, published 1402-02-23
$ws.Range("I9").Value = "1402-02-23 (10)"

# Step 4: fill in the new column M figures (rial, cumulative) for every line item
$ws.Range("M11").Value = 91661180
$ws.Range("M12").Value = -55297002
$ws.Range("M13").Value = 36364178
$ws.Range("M14").Value = -9847885
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 201919
$ws.Range("M17").Value = 26718212
$ws.Range("M18").Value = -5242657
$ws.Range("M19").Value = 403406
$ws.Range("M20").Value = 21878961
$ws.Range("M21").Value = -4312147
$ws.Range("M22").Value = 17566814
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 17566814
$ws.Range("M25").Value = 1039
$ws.Range("M26").Value = 16900000
$ws.Range("M27").Value = 1039

# Step 5: the EPS "read price" algorithm was recomputed for the 1400/09
# column (col I) using the latest capital, same as the latest-capital EPS row
$ws.Range("I25").Value = 591
$ws.Range("I26").Value = 16900000
$ws.Range("K26").Value = 16900000
